$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28: msci_world / benchmark / bonds / 2022 / 0.5 / =1-E28 / 60
$ws.Range("A28").Value = "msci_world"
$ws.Range("B28").Value = "benchmark"
$ws.Range("C28").Value = "bonds"
$ws.Range("D28").Value = 2022
$ws.Range("E28").Value = 0.5
$ws.Range("F28").Formula = "=1-E28"
$ws.Range("G28").Value = 60

# Row 29: msci_world / benchmark / bonds / 2027 / 0.6 / =1-E29 / 75
$ws.Range("A29").Value = "msci_world"
$ws.Range("B29").Value = "benchmark"
$ws.Range("C29").Value = "bonds"
$ws.Range("D29").Value = 2027
$ws.Range("E29").Value = 0.6
$ws.Range("F29").Formula = "=1-E29"
$ws.Range("G29").Value = 75

# Update selection to match the target state
$ws.Range("G30").Select()
